$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.163615822792053
$ws.Range("B1").Value = 2.422133445739746
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.376189708709717
$ws.Range("E1").Value = 1.234842419624329
